# Adds the three new worksheets (04-27, 04-28, 04-29) and fills in the
# previously-empty "04-26 Less und more" sheet, per the "3 neue Blaetter
# hinzugefuegt" commit.

$wb = $excel.ActiveWorkbook

# --- 04-26 Less und more (already existed, was empty) ---------------------
$ws8 = $wb.Worksheets.Item('04-26 Less und more')

$ws8.Range('A1').Value = 'Frage'
$ws8.Range('B1').Value = 'Antwort'
$ws8.Range('A2').Value = 'Zeige die ersten paar Zeilen von spam1.txt an, mit der Möglichkeit (nur) nach unten mit ENTER durchzuscrollen.'
$ws8.Range('B2').Value = 'more spam1.txt'
$ws8.Range('A3').Value = 'Du hast dich vertippt und möchtest deine Eingabe komplett löschen. (Uppercase)'
$ws8.Range('B3').Value = 'STRG + C'
$ws8.Range('A4').Value = 'Zeige die ersten paar Zeilen von spam1.txt an, mit der Möglichkeit mit Pfeiltasten nach unten + oben durchzuscrollen.'
$ws8.Range('B4').Value = 'less spam1.txt'
$ws8.Range('A5').Value = 'Wie beendest du das Programm less? (Uppercase)'
$ws8.Range('B5').Value = 'Q'
$ws8.Columns.Item(1).ColumnWidth = 90.58
$ws8.Columns.Item(2).ColumnWidth = 13.92
$ws8.Activate()
$ws8.Range('B13').Select() | Out-Null

# --- 04-27 Parameter nachschauen (new) ---------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add($null, $last)
$ws9.Name = '04-27 Parameter nachschauen'
$ws9.Range('A1').Value = 'Frage'
$ws9.Range('B1').Value = 'Antwort'
$ws9.Range('A2').Value = 'Rufe die Hilfe zum Befehl "head" mit einem Parameter von "head" auf. Ausführlicher Parameter'
$ws9.Range('B2').Value = 'head --help'
$ws9.Range('A3').Value = 'Aufruf von head mit 4 Zeilen (Kurzschreibweise) von datei text.txt.'
$ws9.Range('B3').Value = 'head -n 4 text.txt'
$ws9.Range('A4').Value = 'Aufruf von head mit 4 Zeilen (Langschreibweise /Verbose) von datei text.txt.'
$ws9.Range('B4').Value = 'head --lines=4 text.txt'
$ws9.Range('A5').Value = 'Ausführliche Dokumentation zum Befehl head anzeigen lassen.'
$ws9.Range('B5').Value = 'man head'
$ws9.Columns.Item(1).ColumnWidth = 86.75
$ws9.Columns.Item(2).ColumnWidth = 19.75
$ws9.Activate()
$ws9.Range('F33').Select() | Out-Null

# --- 04-28 Die Path-Variable (new) ---------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws10 = $wb.Worksheets.Add($null, $last)
$ws10.Name = '04-28 Die Path-Variable'
$ws10.Range('A1').Value = 'Frage'
$ws10.Range('B1').Value = 'Antwort'
$ws10.Range('A2').Value = 'Gib die Variable aus, die alle Pfade zu ausführbaren Befehlen (ohne Pfadangabe) beinhaltet.'
$ws10.Range('B2').Value = 'echo $PATH'
$ws10.Range('A3').Value = 'Gib den Pfad aus, von dem aus der Befehl cat gestartet wird.'
$ws10.Range('B3').Value = 'which cat'
$ws10.Columns.Item(1).ColumnWidth = 83.58
$ws10.Activate()
$ws10.Range('B12').Select() | Out-Null

# --- 04-29 Programme Starten path (new) ---------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws11 = $wb.Worksheets.Add($null, $last)
$ws11.Name = '04-29 Programme Starten path'
$ws11.Range('A1').Value = 'Frage'
$ws11.Range('B1').Value = 'Antwort'
$ws11.Range('A2').Value = 'Ich möchte eine Datei ausführen können. Was sollte ich Ganz oben in die Datei schreiben, wenn die Datei mit bash laufen soll?'
$ws11.Range('B2').Value = '#!/bin/bash'
$ws11.Range('A3').Value = 'Setze die Variable PATH auf den bisherigen Inhalt und füge ~/bin hinzu (erstmal nur temporär)'
$ws11.Range('B3').Value = 'export PATH=$PATH:~/bin'
$ws11.Range('A4').Value = 'Was trennt einzelne Pfade in der PATH-Variable voneinander?'
$ws11.Range('B4').Value = ':'
$ws11.Columns.Item(1).ColumnWidth = 114.58
$ws11.Columns.Item(2).ColumnWidth = 23.42
$ws11.Activate()
$ws11.Range('A3').Select() | Out-Null
